# Auto-generated Excel COM-interop script to apply market-data refresh diff
# Updates columns H-N (price/profit calc columns) on 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1149528.5
$ws.Range("I6").Value = 2296358.5
$ws.Range("J6").Value = 2698.4
$ws.Range("K6").Value = 6889075.5
$ws.Range("L6").Value = 8095.200000000001
$ws.Range("M6").Value = -6888963.5
$ws.Range("N6").Value = -8319.200000000001
$ws.Range("H58").Value = 843053.9
$ws.Range("I58").Value = 1683774.5
$ws.Range("J58").Value = 2333.3333
$ws.Range("K58").Value = 5051323.5
$ws.Range("L58").Value = 6999.999899999999
$ws.Range("M58").Value = -5051173.5
$ws.Range("N58").Value = -7299.999899999999
$ws.Range("H107").Value = 201.57895
$ws.Range("I107").Value = 140.07692
$ws.Range("J107").Value = 334.83334
$ws.Range("K107").Value = 140.07692
$ws.Range("L107").Value = 334.83334
$ws.Range("M107").Value = 1779.92308
$ws.Range("N107").Value = -4174.83334
$ws.Range("H129").Value = 4036.2903
$ws.Range("J129").Value = 950.2857
$ws.Range("L129").Value = 2850.8571
$ws.Range("N129").Value = -12850.8571

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1250.6
$ws.Range("I22").Value = 617.6667
$ws.Range("K22").Value = 617.6667
$ws.Range("M22").Value = -318.6667
$ws.Range("H32").Value = 18803.617
$ws.Range("I32").Value = 3188.896
$ws.Range("J32").Value = 89529.12
$ws.Range("K32").Value = 3188.896
$ws.Range("L32").Value = 89529.12
$ws.Range("M32").Value = -2901.896
$ws.Range("N32").Value = -90103.12
$ws.Range("H61").Value = 2712.36
$ws.Range("I61").Value = 2163.875
$ws.Range("J61").Value = 2970.4707
$ws.Range("K61").Value = 2163.875
$ws.Range("L61").Value = 2970.4707
$ws.Range("M61").Value = -1951.875
$ws.Range("N61").Value = -3394.4707
$ws.Range("H63").Value = 3097.1428
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 3545
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 3545
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -4917
$ws.Range("H66").Value = 3097.1428
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 3545
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 17725
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -24589
$ws.Range("H74").Value = 2380.516
$ws.Range("I74").Value = 1823.3572
$ws.Range("J74").Value = 2839.353
$ws.Range("K74").Value = 1823.3572
$ws.Range("L74").Value = 2839.353
$ws.Range("M74").Value = -949.3571999999999
$ws.Range("N74").Value = -4587.353
$ws.Range("H77").Value = 2380.516
$ws.Range("I77").Value = 1823.3572
$ws.Range("J77").Value = 2839.353
$ws.Range("K77").Value = 9116.786
$ws.Range("L77").Value = 14196.765
$ws.Range("M77").Value = -4748.786
$ws.Range("N77").Value = -22932.765
$ws.Range("H88").Value = 2370.2
$ws.Range("J88").Value = 2462.75
$ws.Range("L88").Value = 2462.75
$ws.Range("N88").Value = -3274.75
$ws.Range("H91").Value = 2370.2
$ws.Range("J91").Value = 2462.75
$ws.Range("L91").Value = 2462.75
$ws.Range("N91").Value = -5270.75
$ws.Range("H102").Value = 60924.65
$ws.Range("I102").Value = 101757.9
$ws.Range("K102").Value = 101757.9
$ws.Range("M102").Value = -100135.9
$ws.Range("H122").Value = 2400.5334
$ws.Range("I122").Value = 2230.5715
$ws.Range("J122").Value = 4780
$ws.Range("K122").Value = 6691.7145
$ws.Range("L122").Value = 14340
$ws.Range("M122").Value = -4241.7145
$ws.Range("N122").Value = -19240
$ws.Range("H136").Value = 2712.36
$ws.Range("I136").Value = 2163.875
$ws.Range("J136").Value = 2970.4707
$ws.Range("K136").Value = 6491.625
$ws.Range("L136").Value = 8911.4121
$ws.Range("M136").Value = -3941.625
$ws.Range("N136").Value = -14011.4121

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66758.414
$ws.Range("I86").Value = 112040.6
$ws.Range("J86").Value = 2069.5715
$ws.Range("K86").Value = 112040.6
$ws.Range("L86").Value = 2069.5715
$ws.Range("M86").Value = -110917.6
$ws.Range("N86").Value = -4315.5715
$ws.Range("H89").Value = 66758.414
$ws.Range("I89").Value = 112040.6
$ws.Range("J89").Value = 2069.5715
$ws.Range("K89").Value = 560203
$ws.Range("L89").Value = 10347.8575
$ws.Range("M89").Value = -554587
$ws.Range("N89").Value = -21579.8575
$ws.Range("H105").Value = 252743.88
$ws.Range("I105").Value = 202592
$ws.Range("J105").Value = 336330.34
$ws.Range("K105").Value = 202592
$ws.Range("L105").Value = 336330.34
$ws.Range("M105").Value = -200845
$ws.Range("N105").Value = -339824.34
$ws.Range("H126").Value = 48000
$ws.Range("J126").Value = 48000
$ws.Range("L126").Value = 48000
$ws.Range("N126").Value = -57880
$ws.Range("H134").Value = 2483.7407
$ws.Range("I134").Value = 2908.1428
$ws.Range("J134").Value = 998.3333
$ws.Range("K134").Value = 8724.428400000001
$ws.Range("L134").Value = 2994.9999
$ws.Range("M134").Value = -6189.428400000001
$ws.Range("N134").Value = -8064.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 695
$ws.Range("I22").Value = 351
$ws.Range("K22").Value = 351
$ws.Range("M22").Value = -1
$ws.Range("H99").Value = 8928.529
$ws.Range("I99").Value = 3854.1428
$ws.Range("K99").Value = 3854.1428
$ws.Range("M99").Value = -2356.1428
$ws.Range("H107").Value = 776.3333
$ws.Range("I107").Value = 904.7
$ws.Range("J107").Value = 615.875
$ws.Range("K107").Value = 904.7
$ws.Range("L107").Value = 615.875
$ws.Range("M107").Value = 1015.3
$ws.Range("N107").Value = -4455.875
$ws.Range("H126").Value = 8928.529
$ws.Range("I126").Value = 3854.1428
$ws.Range("K126").Value = 11562.4284
$ws.Range("M126").Value = -9092.428400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.666668
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 33.666668
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 202.000008
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -428.000008
$ws.Range("H37").Value = 579061.3
$ws.Range("J37").Value = 579061.3
$ws.Range("L37").Value = 1737183.9
$ws.Range("N37").Value = -1737407.9
$ws.Range("H39").Value = 9855.556
$ws.Range("J39").Value = 9855.556
$ws.Range("L39").Value = 29566.668
$ws.Range("N39").Value = -30154.668
$ws.Range("H55").Value = 10650.833
$ws.Range("J55").Value = 12661
$ws.Range("L55").Value = 37983
$ws.Range("N55").Value = -38337
$ws.Range("H64").Value = 1742.4
$ws.Range("I64").Value = 904
$ws.Range("K64").Value = 2712
$ws.Range("M64").Value = -2442
$ws.Range("H67").Value = 1742.4
$ws.Range("I67").Value = 904
$ws.Range("K67").Value = 2712
$ws.Range("M67").Value = -1776
$ws.Range("H68").Value = 1932.7925
$ws.Range("J68").Value = 2650.4194
$ws.Range("L68").Value = 7951.2582
$ws.Range("N68").Value = -9573.2582
$ws.Range("H71").Value = 1932.7925
$ws.Range("J71").Value = 2650.4194
$ws.Range("L71").Value = 23853.7746
$ws.Range("N71").Value = -31965.7746
$ws.Range("H114").Value = 909.2381
$ws.Range("I114").Value = 696
$ws.Range("J114").Value = 994.5333000000001
$ws.Range("K114").Value = 2088
$ws.Range("L114").Value = 2983.5999
$ws.Range("M114").Value = 1166
$ws.Range("N114").Value = -9491.599900000001
$ws.Range("H129").Value = 10418351
$ws.Range("I129").Value = 27778244
$ws.Range("K129").Value = 83334732
$ws.Range("M129").Value = -83329732
$ws.Range("H131").Value = 1542.885
$ws.Range("J131").Value = 1554.4691
$ws.Range("L131").Value = 4663.4073
$ws.Range("N131").Value = -14743.4073

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5349748.5
$ws.Range("I126").Value = 2848
$ws.Range("K126").Value = 8544
$ws.Range("M126").Value = -6074
$ws.Range("H132").Value = 2335.6
$ws.Range("I132").Value = 2202.138
$ws.Range("J132").Value = 2687.4546
$ws.Range("K132").Value = 6606.414
$ws.Range("L132").Value = 8062.3638
$ws.Range("M132").Value = -4076.414
$ws.Range("N132").Value = -13122.3638

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 502750
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -5772
$ws.Range("H132").Value = 8542.5625
$ws.Range("I132").Value = 10964.777
$ws.Range("J132").Value = 5428.2856
$ws.Range("K132").Value = 32894.331
$ws.Range("L132").Value = 16284.8568
$ws.Range("M132").Value = -30364.331
$ws.Range("N132").Value = -21344.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 43980.883
$ws.Range("I107").Value = 10811.454
$ws.Range("J107").Value = 68305.13
$ws.Range("K107").Value = 32434.362
$ws.Range("L107").Value = 204915.39
$ws.Range("M107").Value = -30514.362
$ws.Range("N107").Value = -208755.39
$ws.Range("H136").Value = 727.8043
$ws.Range("I136").Value = 475.25806
$ws.Range("K136").Value = 1425.77418
$ws.Range("M136").Value = 1124.22582
